$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 8.270800000000003
$ws.Range("A3").Value = -21.78090000000001
$ws.Range("C3").Value = -10.996
$ws.Range("C12").Value = -11.25419999999999
$ws.Range("A14").Value = -21.76330000000001
$ws.Range("A21").Value = -20.17159999999999
$ws.Range("A23").Value = -20.28639999999998
$ws.Range("C24").Value = -13.05459999999999
$ws.Range("A25").Value = -21.95409999999999
$ws.Range("B25").Value = 5.866999999999996
$ws.Range("C25").Value = -13.40439999999999
$ws.Range("A26").Value = -21.05239999999996
$ws.Range("B27").Value = 6.064400000000005
$ws.Range("A29").Value = -20.74969999999998
$ws.Range("B31").Value = 5.084199999999998
$ws.Range("B39").Value = 9.858500000000003
$ws.Range("B48").Value = 5.389800000000003
$ws.Range("C50").Value = -13.3089
$ws.Range("B51").Value = 5.674899999999997
$ws.Range("B52").Value = 4.9176
$ws.Range("A53").Value = -21.46320000000001
$ws.Range("C53").Value = -10.2589
$ws.Range("B55").Value = 5.909599999999994
$ws.Range("B56").Value = 4.949599999999998
$ws.Range("A57").Value = -22.092
$ws.Range("B57").Value = 5.267399999999999
$ws.Range("C57").Value = -12.78599999999999
$ws.Range("A59").Value = -22.4796
$ws.Range("C61").Value = -13.03929999999998
$ws.Range("C63").Value = -11.8811
$ws.Range("A69").Value = -21.60239999999999
$ws.Range("C70").Value = -11.567
$ws.Range("B73").Value = 8.3711
$ws.Range("A79").Value = -20.53490000000001
$ws.Range("A83").Value = -21.91269999999999
$ws.Range("C86").Value = -13.2992
$ws.Range("B89").Value = 5.200899999999995
$ws.Range("B90").Value = 5.487300000000001
$ws.Range("A91").Value = -21.29030000000001
$ws.Range("B92").Value = 5.205999999999992
$ws.Range("A93").Value = -20.93089999999998
$ws.Range("C98").Value = -11.63919999999999
$ws.Range("C100").Value = -13.19949999999998
$ws.Range("C102").Value = -13.11310000000001
